# Update release version string across the workbook
$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# A2: "Version: ..."
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

# A6: Recommended citation text
$aboutSheet.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Oaky Creek Coal Mine, Australia, M0085, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# Column S (build_version) rows 2 through 66 on data sheet
for ($r = 2; $r -le 66; $r++) {
    $cell = $dataSheet.Cells.Item($r, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
